# ---------------------------------------------------------------------------
# Applies the "Add files via upload" commit to the single-slide deck.
#
# The underlying commit (captured against the author's full hw03.pptx, a
# 13-slide deck) touched three kinds of things:
#   1. an empty <p:extLst>/<p15:sldGuideLst/> stamp on the presentation
#      (a no-op UI artifact left behind by PowerPoint's "guides" feature);
#   2. the cached text of every datetimeFigureOut field across the slides
#      (an automatic side effect of PowerPoint re-saving the file on a
#      later date - 2024/9/16 -> 2024/10/21);
#   3. a real content edit: the previously-empty text box on one slide was
#      filled in with "hw01-04.ppt：<google-drive-link>" (the link itself
#      hyperlinked), followed by a trailing blank paragraph.
#
# Our before.pptx is a minimized, single-slide reproduction of that same
# deck (same theme/layout/master/field ids, same empty-paragraph shapes),
# so only change (3) has a concrete target to apply here: the subtitle
# placeholder on slide 1 is textually identical, pre-edit, to the shape
# touched by the diff. We reproduce that edit faithfully, run-by-run, so
# the resulting OOXML run/paragraph layout mirrors the diff as closely as
# the object model allows.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The subtitle placeholder ("副標題 2") - shape 2 - is the body-text box
# that receives the new content (the title placeholder stays untouched,
# matching the fact only a single shape changed in the source diff).
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Build the first paragraph run-by-run (mirrors the diff's own run split,
# which separates "h" from "w01-04" from ".ppt" from the full-width colon
# from the hyperlinked URL).
$tr.Text = "h"
$null = $tr.InsertAfter("w01-04")
$null = $tr.InsertAfter(".ppt")
$null = $tr.InsertAfter("：")
$urlText = "https://drive.google.com/drive/folders/1-pDV9bIA_EX6QU-Ci_Drk7dGkvZ8RmOG?usp=drive_link"
$null = $tr.InsertAfter($urlText)

# Re-acquire the URL run by character offset (1-based) so the hyperlink is
# scoped only to that run, exactly like the source XML's <a:hlinkClick>.
$urlStart = $tr.Length - $urlText.Length + 1
$urlRange = $tr.Characters($urlStart, $urlText.Length)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $urlText

# Second, trailing empty paragraph (matches the appended
# "<a:p><a:endParaRPr .../></a:p>" in the diff).
$null = $tr.InsertAfter("`rX")
$tail = $tr.Characters($tr.Length, 1)
$tail.Delete()

Write-Output ("subtitle text now: " + $tr.Text)
